$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update as Text (matches source inlineStr cells),
# guarding against Excel auto-converting numeric-looking strings
# (e.g. "1.000", "0.2660") into numbers, then restoring the default
# "Normal" style so no extra number-format style is introduced.

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '26.489.48'
$cell.Style = 'Normal'

$cell = $ws.Range('E2')
$cell.NumberFormat = '@'
$cell.Value = '  -0.15%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '1.732.04'
$cell.Style = 'Normal'

$cell = $ws.Range('E3')
$cell.NumberFormat = '@'
$cell.Value = '  -0.42%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.Style = 'Normal'

$cell = $ws.Range('E4')
$cell.NumberFormat = '@'
$cell.Value = '  -0.06%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '246.19'
$cell.Style = 'Normal'

$cell = $ws.Range('E5')
$cell.NumberFormat = '@'
$cell.Value = '  +0.25%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E6')
$cell.NumberFormat = '@'
$cell.Value = '  -0.10%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.4872'
$cell.Style = 'Normal'

$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  +1.64%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.2660'
$cell.Style = 'Normal'

$cell = $ws.Range('E8')
$cell.NumberFormat = '@'
$cell.Value = '  -0.80%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.06215'
$cell.Style = 'Normal'

$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  -0.42%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '1.734.54'
$cell.Style = 'Normal'

$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  -0.24%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.07055'
$cell.Style = 'Normal'

$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  -1.21%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E12')
$cell.NumberFormat = '@'
$cell.Value = '  -0.34%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '4.615'
$cell.Style = 'Normal'

$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  +1.94%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '0.6085'
$cell.Style = 'Normal'

$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  -1.39%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '77.28'
$cell.Style = 'Normal'

$cell = $ws.Range('E15')
$cell.NumberFormat = '@'
$cell.Value = '  +0.08%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  -0.06%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '26.485.28'
$cell.Style = 'Normal'

$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  -0.24%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell.Style = 'Normal'

$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  -0.12%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '0.000007227'
$cell.Style = 'Normal'

$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  +4.84%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  -1.90%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '1.956.56'
$cell.Style = 'Normal'

$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  -0.51%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '4.523'
$cell.Style = 'Normal'

$cell = $ws.Range('E22')
$cell.NumberFormat = '@'
$cell.Value = '  -1.55%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '8.738'
$cell.Style = 'Normal'

$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  -1.39%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '5.244'
$cell.Style = 'Normal'

$cell = $ws.Range('E24')
$cell.NumberFormat = '@'
$cell.Value = '  -1.99%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '139.53'
$cell.Style = 'Normal'

$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  +2.85%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  +0.41%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '1.777'
$cell.Style = 'Normal'

$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  -1.55%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '1.407'
$cell.Style = 'Normal'

$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  -1.19%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '108.06'
$cell.Style = 'Normal'

$cell = $ws.Range('E29')
$cell.NumberFormat = '@'
$cell.Value = '  +0.98%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '3.977'
$cell.Style = 'Normal'

$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  -0.08%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '0.08043'
$cell.Style = 'Normal'

$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  +2.39%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '3.683'
$cell.Style = 'Normal'

$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  -1.25%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  +0.21%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B34')
$cell.NumberFormat = '@'
$cell.Value = 'Frax'
$cell.Style = 'Normal'

$cell = $ws.Range('C34')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$cell.Style = 'Normal'

$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.Style = 'Normal'

$cell = $ws.Range('E34')
$cell.NumberFormat = '@'
$cell.Value = '  -0.06%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B35')
$cell.NumberFormat = '@'
$cell.Value = 'HuobiToken'
$cell.Style = 'Normal'

$cell = $ws.Range('C35')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell.Style = 'Normal'

$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '2.615'
$cell.Style = 'Normal'

$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  -0.16%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B36')
$cell.NumberFormat = '@'
$cell.Value = 'ARBITRUM'
$cell.Style = 'Normal'

$cell = $ws.Range('C36')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell.Style = 'Normal'

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '1.009'
$cell.Style = 'Normal'

$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  +1.06%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B37')
$cell.NumberFormat = '@'
$cell.Value = 'ImmutableX'
$cell.Style = 'Normal'

$cell = $ws.Range('C37')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell.Style = 'Normal'

$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '0.6368'
$cell.Style = 'Normal'

$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  +0.22%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B38')
$cell.NumberFormat = '@'
$cell.Value = 'TrustWalletToken'
$cell.Style = 'Normal'

$cell = $ws.Range('C38')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell.Style = 'Normal'

$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.9040'
$cell.Style = 'Normal'

$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  -2.93%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B39')
$cell.NumberFormat = '@'
$cell.Value = 'RenderToken'
$cell.Style = 'Normal'

$cell = $ws.Range('C39')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell.Style = 'Normal'

$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '2.040'
$cell.Style = 'Normal'

$cell = $ws.Range('E39')
$cell.NumberFormat = '@'
$cell.Value = '  +3.55%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B40')
$cell.NumberFormat = '@'
$cell.Value = 'MXToken'
$cell.Style = 'Normal'

$cell = $ws.Range('C40')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell.Style = 'Normal'

$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '2.401'
$cell.Style = 'Normal'

$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  -1.40%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B41')
$cell.NumberFormat = '@'
$cell.Value = 'PaxDollar'
$cell.Style = 'Normal'

$cell = $ws.Range('C41')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell.Style = 'Normal'

$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '1.004'
$cell.Style = 'Normal'

$cell = $ws.Range('E41')
$cell.NumberFormat = '@'
$cell.Value = '  -0.12%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B42')
$cell.NumberFormat = '@'
$cell.Value = 'VeChain'
$cell.Style = 'Normal'

$cell = $ws.Range('C42')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell.Style = 'Normal'

$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.01511'
$cell.Style = 'Normal'

$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  +0.12%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B43')
$cell.NumberFormat = '@'
$cell.Value = 'Quant'
$cell.Style = 'Normal'

$cell = $ws.Range('C43')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$cell.Style = 'Normal'

$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '101.34'
$cell.Style = 'Normal'

$cell = $ws.Range('E43')
$cell.NumberFormat = '@'
$cell.Value = '  -10.25%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B44')
$cell.NumberFormat = '@'
$cell.Value = 'FraxShare'
$cell.Style = 'Normal'

$cell = $ws.Range('C44')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell.Style = 'Normal'

$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '5.437'
$cell.Style = 'Normal'

$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  -5.17%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B45')
$cell.NumberFormat = '@'
$cell.Value = 'TheSandbox'
$cell.Style = 'Normal'

$cell = $ws.Range('C45')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell.Style = 'Normal'

$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.3882'
$cell.Style = 'Normal'

$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  -0.39%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B46')
$cell.NumberFormat = '@'
$cell.Value = 'Aptos'
$cell.Style = 'Normal'

$cell = $ws.Range('C46')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell.Style = 'Normal'

$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '6.958'
$cell.Style = 'Normal'

$cell = $ws.Range('E46')
$cell.NumberFormat = '@'
$cell.Value = '  +3.00%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B47')
$cell.NumberFormat = '@'
$cell.Value = 'Algorand'
$cell.Style = 'Normal'

$cell = $ws.Range('C47')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell.Style = 'Normal'

$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '0.1184'
$cell.Style = 'Normal'

$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  -1.70%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B48')
$cell.NumberFormat = '@'
$cell.Value = 'Cronos'
$cell.Style = 'Normal'

$cell = $ws.Range('C48')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell.Style = 'Normal'

$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '0.05392'
$cell.Style = 'Normal'

$cell = $ws.Range('E48')
$cell.NumberFormat = '@'
$cell.Value = '  +1.13%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '30.54'
$cell.Style = 'Normal'

$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  -0.55%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B50')
$cell.NumberFormat = '@'
$cell.Value = 'EnergySwap'
$cell.Style = 'Normal'

$cell = $ws.Range('C50')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell.Style = 'Normal'

$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '7.796'
$cell.Style = 'Normal'

$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  -1.74%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B51')
$cell.NumberFormat = '@'
$cell.Value = 'NEARProtocol'
$cell.Style = 'Normal'

$cell = $ws.Range('C51')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell.Style = 'Normal'

$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '1.249'
$cell.Style = 'Normal'

$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  -0.48%  '
$cell.Style = 'Normal'

